$d = $word.ActiveDocument

# Locate the title paragraph that currently reads "NUM Test Plan " (bold,
# underlined heading built from two runs: "NUM" and " Test Plan"). We find
# it by scanning paragraphs for the literal text rather than hard-coding an
# index, so the script is resilient to minor structural differences.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "NUM Test Plan*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# Absolute character offsets (relative to the whole document) of the pieces
# of text inside this paragraph:
#   "NUM" " Test Plan" " "
#   [pStart, pStart+3) [pStart+3, pStart+13) [pStart+13, pStart+14)
$numStart = $pStart
$numEnd = $pStart + 3
$secondRunStart = $numEnd
$secondRunEnd = $secondRunStart + 10

# 1) Remove the second run's text (" Test Plan") entirely, leaving the
#    "NUM" run and the trailing " " run, with the (invisible) _GoBack
#    bookmark sitting between the "NUM" run and the trailing space run.
$toDelete = $d.Range($secondRunStart, $secondRunEnd)
$toDelete.Delete()

# 2) Replace the "NUM" run's own text in place (keeps the run, does not
#    merge it away) so it now reads "Test Plan".
$numRange = $d.Range($numStart, $numEnd)
$numRange.Text = "Test Plan"

# 3) Re-home the _GoBack bookmark at the very start of the paragraph so it
#    wraps the "Test Plan" run (matching the target layout where the
#    bookmark precedes the run rather than sitting inside it).
$bmRange = $d.Range($target.Range.Start, $target.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
